# Publish terminology IG 2.0.0
#  - bump Version / Date metadata
#  - insert a new "Properties" sheet (between "Metadata" and "Concepts")
#    describing the CodeSystem's concept properties

$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

# ---------------------------------------------------------------------
# 1. Metadata sheet: bump Version and Date
# ---------------------------------------------------------------------
$metadata.Range("B3").Value = "1.0.1"

# Write the Date cell as literal text (Excel would otherwise auto-convert
# a bare "2025-09-22" into a date serial number). Force text interpretation
# via NumberFormat, then restore the original cell formatting/borders by
# pasting formats back in from the neighbouring (untouched) label cell.
$metadata.Range("B8").NumberFormat = "@"
$metadata.Range("B8").Value = "2025-09-22"
$metadata.Range("A8").Copy()
$metadata.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Insert the new "Properties" sheet between "Metadata" and "Concepts"
#    Copy "Concepts" so the new sheet inherits identical header/body
#    cell styles, then overwrite its contents.
# ---------------------------------------------------------------------
$concepts.Copy($concepts, $null)
$properties = $wb.Worksheets.Item("Concepts (2)")
$properties.Name = "Properties"

# Header row
$properties.Range("A1").Value = "Code"
$properties.Range("B1").Value = "Uri"
$properties.Range("C1").Value = "Description"
$properties.Range("D1").Value = "Type"

# Row 2: status property
$properties.Range("A2").Value = "status"
$properties.Range("B2").Value = "http://hl7.org/fhir/concept-properties#status"
$properties.Range("C2").Value = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$properties.Range("D2").Value = "code"

# Row 3: effectiveDate property (new row - copy formats down from row 2 first)
$properties.Range("A2:D2").Copy()
$properties.Range("A3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$properties.Range("A3").Value = "effectiveDate"
$properties.Range("B3").Value = "http://hl7.org/fhir/concept-properties#effectiveDate"
$properties.Range("C3").Value = "The date at which the concept status was last changed."
$properties.Range("D3").Value = "dateTime"

# ---------------------------------------------------------------------
# 3. Restore original active sheet/selection (Metadata was the tab
#    shown when the workbook was opened - keep it that way).
# ---------------------------------------------------------------------
$metadata.Activate()
$metadata.Range("A1").Select() | Out-Null
